$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell references that receive new text values (numbers/percentages stored as text,
# matching the workbook's existing inline-string convention for columns D and E).
$updates = @{
    "D2" = "296.94"
    "E2" = "-1.40%"
    "D3" = "31.28"
    "E3" = "-0.68%"
    "D4" = "5.074"
    "E4" = "-1.53%"
    "D5" = "0.08015"
    "E5" = "8.89%"
    "D6" = "2.503"
    "E6" = "37.17%"
    "D7" = "7.745"
    "E7" = "-0.61%"
    "D8" = "3.789"
    "E8" = "1.04%"
    "D9" = "0.9308"
    "E9" = "0.17%"
    "D10" = "0.1752"
    "E10" = "3.23%"
    "D11" = "0.07254"
    "E11" = "3.97%"
    "D12" = "0.08998"
    "E12" = "10.83%"
    "D13" = "0.03019"
    "E13" = "-0.46%"
    "D14" = "0.09958"
    "E14" = "0.19%"
    "D15" = "0.001495"
    "E15" = "0.10%"
    "D16" = "0.005914"
    "E16" = "-3.66%"
    "D17" = "3.525"
    "E17" = "1.88%"
    "E18" = "1.25%"
    "D19" = "0.3266"
    "E19" = "-0.29%"
    "D20" = "0.1346"
    "E20" = "1.20%"
    "D21" = "3.390"
    "E21" = "-25.60%"
    "D22" = "0.1618"
    "E22" = "2.37%"
    "D23" = "0.04586"
    "E23" = "-1.22%"
    "D24" = "0.001243"
    "E24" = "2.21%"
    "D25" = "0.004398"
    "E25" = "-7.53%"
    "D26" = "0.0001198"
    "E26" = "-7.61%"
    "D27" = "0.0003428"
    "E27" = "83.17%"
    "D39" = "0.01749"
    "E39" = "1.91%"
    "D40" = "0.04448"
    "E40" = "-1.19%"
    "D41" = "0.006817"
    "E41" = "-4.06%"
    "D42" = "0.1339"
    "E42" = "0.06%"
    "D43" = "0.002141"
    "E43" = "-1.51%"
    "D44" = "0.009580"
    "E44" = "-8.34%"
    "D45" = "0.00006557"
    "E45" = "5.09%"
    "D46" = "0.00000000749"
    "E46" = "-0.19%"
    "D47" = "0.008752"
    "E47" = "-14.30%"
    "E48" = "-55.52%"
    "D49" = "0.00002096"
    "E49" = "-0.19%"
    "D50" = "0.0001996"
    "E50" = "-0.11%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
}
